$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The previous last row (54) had the "date only" number format, marking it
# as the final row. Since row 55 is now the new final row, row 54 goes back
# to the regular datetime number format used by all other data rows, and
# the new row 55 gets the "date only" format instead.
$ws.Cells.Item(54, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"

# Add the new daily row of data.
$ws.Cells.Item(55, 1).Value = 45795
$ws.Cells.Item(55, 1).NumberFormat = "YYYY-MM-DD"
$ws.Cells.Item(55, 2).Value = 230
$ws.Cells.Item(55, 3).Value = 237
$ws.Cells.Item(55, 4).Value = 232
